$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to remain text so numeric-looking values are not
# auto-converted to the Number type (matches original inlineStr/text cells).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '70.525.83'
$ws.Range("E2").Value = '  -3.35%  '
$ws.Range("D3").Value = '3.842.62'
$ws.Range("E3").Value = '  -3.71%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").Value = '592.98'
$ws.Range("E5").Value = '  +0.10%  '
$ws.Range("D6").Value = '167.65'
$ws.Range("E6").Value = '  +5.20%  '
$ws.Range("E7").Value = '  -2.25%  '
$ws.Range("E8").Value = '  +0.19%  '
$ws.Range("E9").Value = '  -0.63%  '
$ws.Range("D10").Value = '0.175'
$ws.Range("E10").Value = '  +3.77%  '
$ws.Range("E11").Value = '  -1.44%  '
$ws.Range("E12").Value = '  +0.01%  '
$ws.Range("D13").Value = '11.34'
$ws.Range("E13").Value = '  +3.60%  '
$ws.Range("D14").Value = '4.461.43'
$ws.Range("E14").Value = '  -3.43%  '
$ws.Range("D15").Value = '21.03'
$ws.Range("E15").Value = '  +3.07%  '
$ws.Range("D16").Value = '3.849.43'
$ws.Range("E16").Value = '  -3.42%  '
$ws.Range("D17").Value = '13.79'
$ws.Range("E17").Value = '  -2.68%  '
$ws.Range("E18").Value = '  -5.89%  '
$ws.Range("E19").Value = '  -2.20%  '
$ws.Range("D20").Value = '70.441.15'
$ws.Range("E20").Value = '  -2.98%  '
$ws.Range("D21").Value = '435.66'
$ws.Range("E21").Value = '  -0.02%  '
$ws.Range("D22").Value = '4.72'
$ws.Range("E22").Value = '  -1.52%  '
$ws.Range("D23").Value = '93.88'
$ws.Range("E23").Value = '  -2.37%  '
$ws.Range("D24").Value = '3.25'
$ws.Range("E24").Value = '  -5.11%  '
$ws.Range("D25").Value = '13.86'
$ws.Range("E25").Value = '  -3.07%  '
$ws.Range("B26").Value = 'RenderToken'
$ws.Range("C26").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D26").Value = '11.19'
$ws.Range("E26").Value = '  -0.41%  '
$ws.Range("B27").Value = 'Toncoin'
$ws.Range("C27").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D27").Value = '4.03'
$ws.Range("E27").Value = '  -10.28%  '
$ws.Range("D28").Value = '5.95'
$ws.Range("E29").Value = '  -1.57%  '
$ws.Range("D30").Value = '35.02'
$ws.Range("E30").Value = '  -3.73%  '
$ws.Range("E31").Value = '  +3.61%  '
$ws.Range("E32").Value = '  -1.68%  '
$ws.Range("D33").Value = '47.92'
$ws.Range("E33").Value = '  -1.36%  '
$ws.Range("E34").Value = '  -4.61%  '
$ws.Range("D35").Value = '69.29'
$ws.Range("E35").Value = '  -0.50%  '
$ws.Range("D36").Value = '0.0₃0977'
$ws.Range("E36").Value = '  +10.52%  '
$ws.Range("D37").Value = '630.82'
$ws.Range("E37").Value = '  -6.82%  '
$ws.Range("D38").Value = '0.423'
$ws.Range("E38").Value = '  -2.99%  '
$ws.Range("B39").Value = 'Dai'
$ws.Range("C39").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D39").Value = '1.00'
$ws.Range("E39").Value = '  -0.15%  '
$ws.Range("B40").Value = 'Kaspa'
$ws.Range("C40").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D40").Value = '0.145'
$ws.Range("E40").Value = '  -1.12%  '
$ws.Range("D41").Value = '0.999'
$ws.Range("E41").Value = '  -0.14%  '
$ws.Range("D42").Value = '3.27'
$ws.Range("E42").Value = '  -2.95%  '
$ws.Range("D43").Value = '3.26'
$ws.Range("E43").Value = '  +23.03%  '
$ws.Range("D44").Value = '0.0468'
$ws.Range("E44").Value = '  -4.02%  '
$ws.Range("D45").Value = '10.01'
$ws.Range("E45").Value = '  -7.63%  '
$ws.Range("D46").Value = '2.71'
$ws.Range("E46").Value = '  +2.60%  '
$ws.Range("E47").Value = '  -4.45%  '
$ws.Range("E48").Value = '  -15.03%  '
$ws.Range("D49").Value = '2.847.77'
$ws.Range("E49").Value = '  +1.31%  '
$ws.Range("D50").Value = '3.23'
$ws.Range("E50").Value = '  -5.38%  '
$ws.Range("E51").Value = '  +0.68%  '

# Restore default styling on column D (clears the temporary text format)
$ws.Range("D2:D51").Style = "Normal"

